$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit: "AF of BC single should be 0 not 100"
# The Attenuation Factor (column G) for rows whose Biochar "type" (column C)
# is "BC _sing" was wrongly recorded as 1 (i.e. 100%). It should be 0.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row   # xlUp = -4162

for ($r = 1; $r -le $lastRow; $r++) {
    $typeVal = $ws.Cells.Item($r, 3).Value2
    if ($typeVal -eq "BC _sing") {
        $afCell = $ws.Cells.Item($r, 7)
        if ($afCell.Value2 -eq 1) {
            $afCell.Value = 0
        }
    }
}

# Minor floating-point refresh of the BC _S _mix attenuation factor in row 9
# (same value, last-bit precision refresh from the upstream recalculation).
$ws.Cells.Item(9, 7).Value = 0.41997890083802097

# Reset the view: scroll back to top-left and select G1 (matches the saved view state).
$ws.Range("G1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
